# Investment Excel - 3rd Sheet added
#
# Fills in the answer values on the "Table - 2.1" (average funding by type)
# and "Table-3.1" (top English speaking countries) sheets, attaches a
# comment with summary stats to C9 on "Table - 2.1", and updates the
# selections / active sheet to match the author's final view.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Table - 2.1" (2nd sheet): average funding amounts per type
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("C5").Value = 11748949.130000001
$ws2.Range("C6").Value = 958694.47
$ws2.Range("C7").Value = 719818
$ws2.Range("C8").Value = 73308593.030000001
$ws2.Range("C9").Value = "Venture"

$comment = $ws2.Range("C9").AddComment("Bhushan, Shashi :`r`n`r`nAverage: 11 Million`r`nmin 0.000000e+00 25% 1.600902e+06 50% 5.000000e+06 75% 1.200000e+07 max 1.760000e+10`r`n")

# ---------------------------------------------------------------------
# Sheet "Table-3.1" (3rd sheet): top 3 English speaking countries
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Range("C5").Value = "USA"
$ws3.Range("C6").Value = "GBR"
$ws3.Range("C7").Value = "CHN"

# ---------------------------------------------------------------------
# Selections / active sheet: "Table-3.1" ends up the active tab, with
# "Table - 2.1" left selected at F9 and "Table-3.1" selected at C10.
# ---------------------------------------------------------------------
$ws2.Activate() | Out-Null
$ws2.Range("F9").Select() | Out-Null

$ws3.Activate() | Out-Null
$ws3.Range("C10").Select() | Out-Null
